# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Row 3: Intel(R) Wireless-AC 9260 160MHz - 23.40.0.4
$ws.Range("C3").Value = 366

# Row 5: Intel(R) Dual Band Wireless-AC 8265 - 20.70.25.2
$ws.Range("C5").Value = 1036
$ws.Range("D5").Value = 98.59999999999999

# Row 7: Totals
$ws.Range("C7").Value = 2511
